# Update sheet1 data to reflect new TPM values and drop the MuSCs-sending rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (old "MuSCs" sending-cluster rows) entirely.
$ws.Rows("8:10").Delete()

# Row 2: FAPs -> Pgf/Nrp1 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Pgf"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.736173
$ws.Range("H2").Value = 11.208519
$ws.Range("I2").Value = 0.8038593646893297
$ws.Range("J2").Value = 0.8038593646893297
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 413.380298276484
$ws.Range("R2").Value = 3720.422684488356
$ws.Range("S2").Value = 0.4402270638133865
$ws.Range("T2").Value = 0.4402270638133865

# Row 3: FAPs -> Pgf/Nrp1 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Pgf"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.736173
$ws.Range("H3").Value = 11.208519
$ws.Range("I3").Value = 0.8038593646893297
$ws.Range("J3").Value = 0.8038593646893297
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 237.854593088035
$ws.Range("R3").Value = 2140.691337792315
$ws.Range("S3").Value = 0.2533019342388677
$ws.Range("T3").Value = 0.2533019342388677

# Row 4: FAPs -> Pgf/Nrp1 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Pgf"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.736173
$ws.Range("H4").Value = 11.208519
$ws.Range("I4").Value = 0.8038593646893297
$ws.Range("J4").Value = 0.8038593646893297
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 103.601990014684
$ws.Range("R4").Value = 932.417910132156
$ws.Range("S4").Value = 0.1103303666370755
$ws.Range("T4").Value = 0.1103303666370755

# Row 5: MuSCs -> Pgf/Nrp1 -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pgf"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9116213333333333
$ws.Range("H5").Value = 2.734864
$ws.Range("I5").Value = 0.1961406353106703
$ws.Range("J5").Value = 0.1961406353106703
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 100.8642529905706
$ws.Range("R5").Value = 907.7782769151358
$ws.Range("S5").Value = 0.1074148287252699
$ws.Range("T5").Value = 0.1074148287252699

# Row 6: MuSCs -> Pgf/Nrp1 -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pgf"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9116213333333333
$ws.Range("H6").Value = 2.734864
$ws.Range("I6").Value = 0.1961406353106703
$ws.Range("J6").Value = 0.1961406353106703
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 58.03621012473778
$ws.Range("R6").Value = 522.32589112264
$ws.Range("S6").Value = 0.0618053411945188
$ws.Range("T6").Value = 0.0618053411945188

# Row 7: MuSCs -> Pgf/Nrp1 -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pgf"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9116213333333333
$ws.Range("H7").Value = 2.734864
$ws.Range("I7").Value = 0.1961406353106703
$ws.Range("J7").Value = 0.1961406353106703
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 25.27875028088177
$ws.Range("R7").Value = 227.508752527936
$ws.Range("S7").Value = 0.0269204653908816
$ws.Range("T7").Value = 0.0269204653908816
